$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '64.269.79'
$ws.Cells.Item(2, 5).Value = '  +8.72%  '
$ws.Cells.Item(3, 4).Value = '3.154.27'
$ws.Cells.Item(3, 5).Value = '  +6.13%  '
$ws.Cells.Item(4, 5).Value = '  +0.20%  '
$ws.Cells.Item(5, 4).Value = '590.86'
$ws.Cells.Item(5, 5).Value = '  +4.50%  '
$ws.Cells.Item(6, 4).Value = '148.13'
$ws.Cells.Item(6, 5).Value = '  +8.31%  '
$ws.Cells.Item(7, 5).Value = '  +0.10%  '
$ws.Cells.Item(8, 4).Value = '3.146.23'
$ws.Cells.Item(8, 5).Value = '  +6.01%  '
$ws.Cells.Item(9, 5).Value = '  +3.86%  '
$ws.Cells.Item(10, 4).Value = '0.158'
$ws.Cells.Item(10, 5).Value = '  +19.56%  '
$ws.Cells.Item(11, 4).Value = '5.80'
$ws.Cells.Item(11, 5).Value = '  +10.48%  '
$ws.Cells.Item(12, 4).Value = '0.472'
$ws.Cells.Item(12, 5).Value = '  +5.43%  '
$ws.Cells.Item(13, 5).Value = '  +11.11%  '
$ws.Cells.Item(14, 4).Value = '35.93'
$ws.Cells.Item(14, 5).Value = '  +7.23%  '
$ws.Cells.Item(15, 4).Value = '0.123'
$ws.Cells.Item(15, 5).Value = '  +1.19%  '
$ws.Cells.Item(16, 4).Value = '3.683.74'
$ws.Cells.Item(16, 5).Value = '  +6.46%  '
$ws.Cells.Item(17, 4).Value = '64.219.21'
$ws.Cells.Item(17, 5).Value = '  +8.67%  '
$ws.Cells.Item(18, 5).Value = '  +2.16%  '
$ws.Cells.Item(19, 4).Value = '3.154.76'
$ws.Cells.Item(19, 5).Value = '  +6.26%  '
$ws.Cells.Item(20, 4).Value = '476.06'
$ws.Cells.Item(20, 5).Value = '  +9.45%  '
$ws.Cells.Item(21, 4).Value = '14.28'
$ws.Cells.Item(21, 5).Value = '  +5.11%  '
$ws.Cells.Item(22, 4).Value = '0.735'
$ws.Cells.Item(22, 5).Value = '  +2.13%  '
$ws.Cells.Item(23, 4).Value = '7.67'
$ws.Cells.Item(23, 5).Value = '  +9.64%  '
$ws.Cells.Item(24, 4).Value = '13.40'
$ws.Cells.Item(24, 5).Value = '  +2.49%  '
$ws.Cells.Item(25, 4).Value = '82.89'
$ws.Cells.Item(25, 5).Value = '  +3.90%  '
$ws.Cells.Item(26, 5).Value = '  +0.05%  '
$ws.Cells.Item(27, 5).Value = '  +13.46%  '
$ws.Cells.Item(28, 2).Value = 'PancakeSwap'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(28, 4).Value = '2.71'
$ws.Cells.Item(28, 5).Value = '  +6.51%  '
$ws.Cells.Item(29, 2).Value = 'ImmutableX'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(29, 4).Value = '2.24'
$ws.Cells.Item(29, 5).Value = '  +1.16%  '
$ws.Cells.Item(30, 5).Value = '  +0.20%  '
$ws.Cells.Item(31, 4).Value = '6.91'
$ws.Cells.Item(31, 5).Value = '  +12.20%  '
$ws.Cells.Item(32, 4).Value = '27.31'
$ws.Cells.Item(32, 5).Value = '  +6.55%  '
$ws.Cells.Item(33, 5).Value = '  +6.43%  '
$ws.Cells.Item(34, 4).Value = '0.0₃0890'
$ws.Cells.Item(34, 5).Value = '  +16.68%  '
$ws.Cells.Item(35, 4).Value = '2.43'
$ws.Cells.Item(35, 5).Value = '  +18.51%  '
$ws.Cells.Item(36, 5).Value = '  +7.76%  '
$ws.Cells.Item(37, 4).Value = '3.44'
$ws.Cells.Item(37, 5).Value = '  +24.36%  '
$ws.Cells.Item(38, 4).Value = '6.16'
$ws.Cells.Item(38, 5).Value = '  +4.98%  '
$ws.Cells.Item(39, 4).Value = '50.87'
$ws.Cells.Item(39, 5).Value = '  +5.06%  '
$ws.Cells.Item(40, 4).Value = '452.09'
$ws.Cells.Item(40, 5).Value = '  +14.30%  '
$ws.Cells.Item(41, 4).Value = '8.78'
$ws.Cells.Item(41, 5).Value = '  +1.02%  '
$ws.Cells.Item(42, 4).Value = '0.0375'
$ws.Cells.Item(42, 5).Value = '  +7.18%  '
$ws.Cells.Item(43, 4).Value = '2.943.93'
$ws.Cells.Item(43, 5).Value = '  +8.29%  '
$ws.Cells.Item(44, 4).Value = '0.284'
$ws.Cells.Item(44, 5).Value = '  +14.12%  '
$ws.Cells.Item(45, 5).Value = '  +6.39%  '
$ws.Cells.Item(46, 5).Value = '  +13.20%  '
$ws.Cells.Item(47, 4).Value = '35.71'
$ws.Cells.Item(47, 5).Value = '  +4.14%  '
$ws.Cells.Item(48, 5).Value = '  +0.03%  '
$ws.Cells.Item(49, 4).Value = '123.61'
$ws.Cells.Item(49, 5).Value = '  +0.98%  '
$ws.Cells.Item(50, 5).Value = '  +2.72%  '
$ws.Cells.Item(51, 4).Value = '25.08'
$ws.Cells.Item(51, 5).Value = '  +8.63%  '
